$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Kaeng Suea Ten" reservoir as row 7 of the data table
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Kaeng Suea Ten"
$ws.Range("C7").Value = 48
$ws.Range("E7").Value = 0.177
$ws.Range("D7").Formula = "=C7*1000000/E7/1000/3600*86400*30.25/1000000"
$ws.Range("F7").Value = 1175

# Row heights: header row slightly shorter, data rows back to default (auto) height
$ws.Rows(1).RowHeight = 30
$ws.Rows(2).AutoFit()
$ws.Rows(3).AutoFit()

# Restore the active selection to where the user last left off
$ws.Range("E23").Select() | Out-Null
